# Weekly price update: a new weekly record is inserted at row 18 (this
# dataset is ordered with the most recent weeks near the top of the
# "recent" block starting at row 18), pushing all subsequent rows down by
# one. The table therefore grows from 76 to 77 data-bearing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; Excel shifts rows 18-76 down to
# 19-77 automatically, carrying their values/formatting with them.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Cells.Item(18, 1).Value  = 9
$ws.Cells.Item(18, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value  = "Metropolitana"
$ws.Cells.Item(18, 4).Value  = 44648
$ws.Cells.Item(18, 5).Value  = 13
$ws.Cells.Item(18, 6).Value  = 100114002
$ws.Cells.Item(18, 7).Value  = "Camote"
$ws.Cells.Item(18, 8).Value  = "Sin especificar"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 1600
$ws.Cells.Item(18, 11).Value = 11000
$ws.Cells.Item(18, 12).Value = 12000
$ws.Cells.Item(18, 13).Value = 11500
$ws.Cells.Item(18, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 639
$ws.Cells.Item(18, 17).Value = 18
$ws.Cells.Item(18, 18).Value = "Hortaliza"
